# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values for rows 2-37 (column G), replacing the old Strike# totals.
$kValues = @{
    2  = 4
    3  = 3
    4  = 6
    5  = 3
    6  = 0
    7  = 7
    8  = 4
    9  = 2
    10 = 3
    11 = 5
    12 = 2
    13 = 4
    14 = 3
    15 = 4
    16 = 6
    17 = 2
    18 = 3
    19 = 4
    20 = 2
    21 = 8
    22 = 2
    23 = 5
    24 = 2
    25 = 6
    26 = 5
    27 = 4
    28 = 5
    29 = 3
    30 = 4
    31 = 9
    32 = 3
    33 = 3
    34 = 4
    35 = 2
    36 = 3
    37 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
